$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "only_schema" column (column C) - the whole column shifts left
$ws.Columns.Item(3).Delete()

# Set the selection to match the post-delete state (Excel leaves the entire
# deleted column's former position selected as the new column at that index)
$ws.Range("C1:C1048576").Select()
